# cryptos.xlsx refresh -- updated Price (D) / Volume(1h) (E) columns, and
# re-ranked Toncoin above BinanceUSD (rows 24-25), per the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store numeric-looking prices (e.g. '0.257')
# as literal text, matching the workbook's existing inlineStr/text cells instead
# of letting COM auto-coerce them into floating-point numbers.
$q = [char]39

$ws.Range('D2').Value = '25.816.92'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.637.57'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = $q + '0.257'
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').Value = $q + '19.79'
$ws.Range('E10').Value = '  -2.50%  '
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '1.863.05'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '1.634.07'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').Value = '0.0₃0769'
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = $q + '63.29'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '25.835.41'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = $q + '4.48'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = $q + '192.87'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').Value = $q + '9.99'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = $q + '6.36'
$ws.Range('E23').Value = '  +1.43%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = $q + '1.82'
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = $q + '1.00'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').Value = $q + '142.42'
$ws.Range('E26').Value = '  +2.33%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = $q + '6.97'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').Value = $q + '15.56'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  -0.83%  '
$ws.Range('D31').Value = $q + '0.0496'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').Value = $q + '3.35'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('D36').Value = $q + '0.907'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').Value = '1.131.35'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  -2.04%  '
$ws.Range('D39').Value = $q + '0.545'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = $q + '5.55'
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').Value = $q + '100.45'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').Value = '1.772.58'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('D47').Value = $q + '55.39'
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('E48').Value = '  -1.59%  '
$ws.Range('E49').Value = '  -0.34%  '
$ws.Range('E50').Value = '  +3.46%  '
$ws.Range('D51').Value = $q + '2.32'
$ws.Range('E51').Value = '  +2.84%  '
